$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("A1").Value = "Danh mục lương"
$ws.Range("B1").Value = 17

$ws.Range("A2").Value = "Ngày công"
$ws.Range("B2").Value = 23.5

$ws.Range("A3").Value = "Phụ cấp"
$ws.Range("B3").Value = 822500

$ws.Range("A4").Value = "Lương cơ bản tại CẦN THƠ"
$ws.Range("B4").ClearContents()

$ws.Range("A5").Value = "Chiết khấu sale chính tại CẦN THƠ"
$ws.Range("B5").Value = 0

$ws.Range("A6").Value = "Chiết khấu sale phụ tại CẦN THƠ"
$ws.Range("B6").Value = 0

$ws.Range("A7").Value = "Đơn 1 bác sĩ tại CẦN THƠ"
$ws.Range("B7").Value = 0

$ws.Range("A8").Value = "Đơn 2 bác sĩ tại CẦN THƠ"
$ws.Range("B8").Value = 0

$ws.Range("A9").Value = "Công phụ phẫu 1 tại CẦN THƠ"
$ws.Range("B9").Value = 0

$ws.Range("A10").Value = "Công phụ phẫu 2 tại CẦN THƠ"
$ws.Range("B10").Value = 0

$ws.Range("A11").Value = "Ứng lương tại CẦN THƠ"
$ws.Range("B11").Value = -0

$ws.Range("A12").Value = "Lương cơ bản tại LONG XUYÊN"
$ws.Range("B12").Value = 2517857.142857143

$ws.Range("A13").Value = "Chiết khấu sale chính tại LONG XUYÊN"
$ws.Range("B13").Value = 0

$ws.Range("A14").Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws.Range("B14").Value = 0

$ws.Range("A15").Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws.Range("B15").Value = 0

$ws.Range("A16").Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws.Range("B16").Value = 0

$ws.Range("A17").Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws.Range("B17").Value = 250000

$ws.Range("A18").Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws.Range("B18").Value = 0

$ws.Range("A19").Value = "Ứng lương tại LONG XUYÊN"
$ws.Range("B19").Value = -1500000

$ws.Range("A20").Value = "Lương cơ bản tại SÓC TRĂNG"
$ws.Range("B20").ClearContents()

$ws.Range("A21").Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws.Range("B21").Value = 0

$ws.Range("A22").Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws.Range("B22").Value = 0

$ws.Range("A23").Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws.Range("B23").Value = 0

$ws.Range("A24").Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws.Range("B24").Value = 0

$ws.Range("A25").Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws.Range("B25").Value = 0

$ws.Range("A26").Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws.Range("B26").Value = 0

$ws.Range("A27").Value = "Ứng lương tại SÓC TRĂNG"
$ws.Range("B27").Value = -0

$ws.Range("A28").Value = "Tổng lương tại CẦN THƠ"
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = "Tổng lương tại LONG XUYÊN"
$ws.Range("B29").Value = 2090357.142857143

$ws.Range("A30").Value = "Tổng lương tại SÓC TRĂNG"
$ws.Range("B30").Value = 0

$ws.Range("A31").Value = "Tổng lương"
$ws.Range("B31").Value = 2090357.142857143

$ws.Range("A32:B39").ClearContents()
